$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# --- Row 16: fill the previously-empty columns B..K with the literal text "nan" ---
$cols16 = 2..11  # B..K
foreach ($col in $cols16) {
    $cell = $ws.Cells.Item(16, $col)
    $cell.Value = "nan"
}

# --- Row 17: brand-new service-log row ---
# A17 = "22" (card number, stored as text like the rest of column A)
$cellA17 = $ws.Cells.Item(17, 1)
$cellA17.NumberFormat = "@"
$cellA17.Value = "22"
$cellA17.Style = "Normal"

# B17..K17 stay blank (no data reported for this event)

# L17 = date of the event
$ws.Cells.Item(17, 12).Value = "25\1\2025"

# M17 stays blank (no "Event" hours entry for this row)

# N17 = correction/notes text
$ws.Cells.Item(17, 14).Value = "تم سن الفلاتس والسليندر وتغيير الجرائد الخلفيه (1_5_8)"

# O17 = serviced by
$ws.Cells.Item(17, 15).Value = "الخبير"
